# Update "想去人数" (want-to-go count) figures across sheets, per the
# upstream data refresh (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 4430
$ws.Range("F6").Value = 429
$ws.Range("F7").Value = 3542
$ws.Range("F8").Value = 1009
$ws.Range("F11").Value = 322
$ws.Range("F12").Value = 317
$ws.Range("F13").Value = 2406
$ws.Range("F18").Value = 528
$ws.Range("F20").Value = 58
$ws.Range("F21").Value = 9974
$ws.Range("F22").Value = 5949
$ws.Range("F23").Value = 381
$ws.Range("F25").Value = 821
$ws.Range("F27").Value = 838
$ws.Range("F28").Value = 3528
$ws.Range("F31").Value = 456
$ws.Range("F33").Value = 239
$ws.Range("F34").Value = 235
$ws.Range("F35").Value = 223
$ws.Range("F36").Value = 4808
$ws.Range("F38").Value = 1080
$ws.Range("F40").Value = 20
$ws.Range("F41").Value = 57

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 3523

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8694
$ws.Range("F3").Value = 418
$ws.Range("F4").Value = 1555

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8694
$ws.Range("F4").Value = 1555
$ws.Range("F6").Value = 4430
$ws.Range("F8").Value = 429
$ws.Range("F9").Value = 3542
$ws.Range("F10").Value = 1009
$ws.Range("F13").Value = 2406
$ws.Range("F22").Value = 528
$ws.Range("F24").Value = 58
$ws.Range("F25").Value = 9975
$ws.Range("F26").Value = 3523
$ws.Range("F28").Value = 381
$ws.Range("F30").Value = 821
$ws.Range("F32").Value = 838
$ws.Range("F33").Value = 3528
$ws.Range("F36").Value = 456
$ws.Range("F38").Value = 235
$ws.Range("F40").Value = 223
$ws.Range("F41").Value = 4808
$ws.Range("F42").Value = 1080
$ws.Range("F44").Value = 57
